# feat: Implement map visualization
#
# Adds a new "Tiles View" worksheet (after "Tile Types") that maps each
# tile-type row to its in-map prefab "view" name and the resource path of
# the prefab used to render it on the map.

$wb = $excel.ActiveWorkbook

# Append the new sheet at the end of the workbook (after the last existing
# sheet, "Tile Types"), matching the tab order in the target workbook.
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Tiles View"

# Header row.
$ws.Cells.Item(1, 1).Value = "IDS"
$ws.Cells.Item(1, 2).Value = "Architecture ID"
$ws.Cells.Item(1, 3).Value = "Prefab resource path"

# One row per tile type defined in the "Tile Types" sheet: the map-view
# display name, the tile type it corresponds to, and the prefab resource
# path used to render that tile on the map.
$rows = @(
    @("Road view",          "Road",          "Prefabs/Map/YellowTile"),
    @("Jail Walls view",    "Jail Walls",    "Prefabs/Map/RedTile"),
    @("Jail Habitat view",  "Jail Habitat",  "Prefabs/Map/CyanTile"),
    @("Structure",          "Power Supply",  "Prefabs/Map/GreenTile"),
    @("Humans Entry View",  "Human Entry",   "Prefabs/Map/BlueTile"),
    @("Humans Exit View",   "Human Exit",    "Prefabs/Map/PinkTile"),
    @("Grass View",         "Grass",         "Prefabs/Map/OrangeTile")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
